$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$ws.Range("C145").Value = 1.171881684295534
$ws.Range("C147").Value = 1.217823389384934
$ws.Range("C149").Value = 1.228879192361553
$ws.Range("C150").Value = 1.226892353205352
$ws.Range("C151").Value = 1.249117742245347
$ws.Range("C152").Value = 1.23607283366597
$ws.Range("C153").Value = 1.224510215525612
$ws.Range("C156").Value = 1.23730038985183
$ws.Range("C157").Value = 1.228944539426958
$ws.Range("C158").Value = 1.227174753504499
$ws.Range("C160").Value = 1.22466122937041
$ws.Range("C163").Value = 1.230561082727156
$ws.Range("C167").Value = 1.259691051498418
$ws.Range("C175").Value = 1.251640714470686
$ws.Range("C176").Value = 1.240701803075182
$ws.Range("C178").Value = 1.239229319877219
$ws.Range("C179").Value = 1.239385813920769
$ws.Range("C184").Value = 1.230371637033066
$ws.Range("C185").Value = 1.229065773288799
$ws.Range("C186").Value = 1.228182306305312
$ws.Range("C187").Value = 1.214760197372009
$ws.Range("C188").Value = 1.214737649731441
$ws.Range("C189").Value = 1.216393863917958
$ws.Range("C190").Value = 1.195224373986969
$ws.Range("C191").Value = 1.206114938435964
$ws.Range("C192").Value = 1.215058924347624
$ws.Range("C193").Value = 1.213940947482899
$ws.Range("C194").Value = 1.21513958041862
$ws.Range("C195").Value = 1.192994489511771
$ws.Range("C196").Value = 1.192555829941442
$ws.Range("C197").Value = 1.192073467763204
